# Apply cell updates to match the refreshed cryptos snapshot.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.606.60'
$ws.Range("E2").Value = '  -2.53%  '
$ws.Range("D3").Value = '1.658.70'
$ws.Range("E3").Value = '  -4.25%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '215.10'
$ws.Range("D5").NumberFormat = "General"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.99%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.509'
$ws.Range("D6").NumberFormat = "General"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -2.77%  '
$ws.Range("E7").Value = '  +0.02%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '24.11'
$ws.Range("D8").NumberFormat = "General"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.30%  '
$ws.Range("E9").Value = '  -2.63%  '
$ws.Range("E10").Value = '  -2.80%  '
$ws.Range("E11").Value = '  -1.79%  '
$ws.Range("D12").Value = '1.893.99'
$ws.Range("E12").Value = '  -4.18%  '
$ws.Range("D13").Value = '1.655.95'
$ws.Range("E13").Value = '  -4.42%  '
$ws.Range("E14").Value = '  -2.91%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.566'
$ws.Range("D15").NumberFormat = "General"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.14%  '
$ws.Range("E16").Value = '  -2.87%  '
$ws.Range("D17").Value = '27.592.01'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '240.89'
$ws.Range("D18").NumberFormat = "General"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.17%  '
$ws.Range("E19").Value = '  -3.51%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.60'
$ws.Range("D20").NumberFormat = "General"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -4.61%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.47'
$ws.Range("D22").NumberFormat = "General"
$ws.Range("D22").Style = "Normal"
$ws.Range("E23").Value = '  -3.99%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.05'
$ws.Range("D24").NumberFormat = "General"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -2.72%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '146.06'
$ws.Range("D25").NumberFormat = "General"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -2.24%  '
$ws.Range("E26").Value = '  -4.68%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '16.29'
$ws.Range("D27").NumberFormat = "General"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -2.21%  '
$ws.Range("E28").Value = '  +0.02%  '
$ws.Range("E29").Value = '  -2.57%  '
$ws.Range("E30").Value = '  +0.37%  '
$ws.Range("E31").Value = '  -2.97%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.33'
$ws.Range("D32").NumberFormat = "General"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -2.82%  '
$ws.Range("D33").Value = '1.456.90'
$ws.Range("E33").Value = '  -2.16%  '
$ws.Range("E34").Value = '  -5.11%  '
$ws.Range("E35").Value = '  -5.04%  '
$ws.Range("E36").Value = '  -1.10%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.925'
$ws.Range("D37").NumberFormat = "General"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -5.37%  '
$ws.Range("E38").Value = '  -2.53%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.573'
$ws.Range("D39").NumberFormat = "General"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -5.05%  '
$ws.Range("E40").Value = '  -0.73%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.02'
$ws.Range("D41").NumberFormat = "General"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -4.28%  '
$ws.Range("E42").Value = '  -0.03%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.42'
$ws.Range("D43").NumberFormat = "General"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -4.26%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.792'
$ws.Range("D45").NumberFormat = "General"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.51%  '
$ws.Range("D46").Value = '1.801.67'
$ws.Range("E46").Value = '  -4.12%  '
$ws.Range("E47").Value = '  -0.85%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '88.61'
$ws.Range("D48").NumberFormat = "General"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -2.50%  '
$ws.Range("E49").Value = '  -6.09%  '
$ws.Range("E50").Value = '  -1.92%  '
$ws.Range("E51").Value = '  -5.23%  '
